$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 107, pushing current rows 107-109 down to 109-111.
$ws.Range("A107:R108").EntireRow.Insert()

# New row 107 data
$ws.Range("A107").Value = 6
$ws.Range("B107").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C107").Value = 'Metropolitana'
$ws.Range("D107").Value = 44448
$ws.Range("D107").NumberFormat = $ws.Range("D106").NumberFormat
$ws.Range("E107").Value = 13
$ws.Range("F107").Value = 100112026
$ws.Range("G107").Value = 'Haba'
$ws.Range("H107").Value = 'Sin especificar'
$ws.Range("I107").Value = 'Primera'
$ws.Range("J107").Value = 250
$ws.Range("K107").Value = 12000
$ws.Range("L107").Value = 13000
$ws.Range("M107").Value = 12400
$ws.Range("N107").Value = '$/saco 25 kilos'
$ws.Range("O107").Value = 'Región Metropolitana'
$ws.Range("P107").Value = 496
$ws.Range("Q107").Value = 25
$ws.Range("R107").Value = 'Hortaliza'

# New row 108 data
$ws.Range("A108").Value = 6
$ws.Range("B108").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C108").Value = 'Metropolitana'
$ws.Range("D108").Value = 44448
$ws.Range("D108").NumberFormat = $ws.Range("D106").NumberFormat
$ws.Range("E108").Value = 13
$ws.Range("F108").Value = 100112026
$ws.Range("G108").Value = 'Haba'
$ws.Range("H108").Value = 'Sin especificar'
$ws.Range("I108").Value = 'Primera'
$ws.Range("J108").Value = 700
$ws.Range("K108").Value = 13000
$ws.Range("L108").Value = 14000
$ws.Range("M108").Value = 13571
$ws.Range("N108").Value = '$/saco 25 kilos'
$ws.Range("O108").Value = 'Región de Coquimbo'
$ws.Range("P108").Value = 543
$ws.Range("Q108").Value = 25
$ws.Range("R108").Value = 'Hortaliza'
